{"js": "// Fill in the student-info table on the cover page with the submitted\n// values (Roll No., Name, Class, Batch, Date of Experiment). The\n// \"Date/Time of Submission\" label itself is left blank (no value was\n// filled in for it in the source edit).\nconst body = context.document.body;\n\nconst fills = [\n  { label: \"Roll No. :\", value: \" C146\" },\n  { label: \"Name:\", value: \" Manan Gandhi\" },\n  { label: \"Class :\", value: \" D\" },\n  { label: \"Batch :\", value: \" D1\" },\n  { label: \"Date of Experiment :\", value: \" 24/07/2024\" }\n];\n\nfor (const fill of fills) {\n  const results = body.search(fill.label, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length > 0) {\n    results.items[0].insertText(fill.value, \"After\");\n    await context.sync();\n  }\n}\n", "ps1": "# Fill in the student-info table on the cover page with the submitted\n# values (Roll No., Name, Class, Batch, Date of Experiment). The\n# \"Date/Time of Submission\" label itself is left blank (no value was\n# filled in for it in the source edit).\n$d = $word.ActiveDocument\n\nfunction Fill-After([string]$label, [string]$value) {\n    $rng = $d.Content\n    $found = $rng.Find.Execute($label)\n    if ($found) {\n        $rng.Collapse(0)\n        $rng.InsertAfter($value)\n    }\n}\n\nFill-After \"Roll No. :\" \" C146\"\nFill-After \"Name:\" \" Manan Gandhi\"\nFill-After \"Class :\" \" D\"\nFill-After \"Batch :\" \" D1\"\nFill-After \"Date of Experiment :\" \" 24/07/2024\"\n"}
